$d = $word.ActiveDocument

$newText = "Perioadele campaniei din Orion: 16-25 ianuarie, 14-23 februarie, 14-24 martie"
$marker = "Perioadele campaniei"

# Collect the paragraphs that still contain the old campaign-dates sentence
# (there are four near-duplicate occurrences throughout the guide).
$targets = New-Object System.Collections.ArrayList
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*$marker*") {
        [void]$targets.Add($i)
    }
}

foreach ($i in $targets) {
    $p = $d.Paragraphs($i)
    $r = $p.Range
    # Exclude the trailing paragraph mark so we only touch the run(s) that
    # hold the sentence itself.
    $r.End = $r.End - 1
    $r.Delete()
    # InsertAfter on a collapsed/emptied range creates a brand-new run with
    # no rPr (no inherited character formatting), which is what the target
    # markup needs.
    $r.InsertAfter($newText)
}
